$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B11 currently holds the text "R40" (shared string). The edit changes it to
# the text "1". A leading apostrophe forces Excel to store it as text rather
# than coercing the numeric-looking string to a number, keeping the cell's
# existing type (shared string) consistent with the rest of the table.
$ws.Range("B11").Value = "'1"
